$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.621.93'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '3.133.01'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.78'
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.52'
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.130.89'
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("E10").Value = '  +4.25%  '
$ws.Range("E11").Value = '  -1.75%  '
$ws.Range("E12").Value = '  -3.19%  '
$ws.Range("E13").Value = '  -1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.99'
$ws.Range("E14").Value = '  +2.79%  '
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").Value = '3.654.13'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.142.54'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '63.526.08'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '462.50'
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.92'
$ws.Range("E24").Value = '  -3.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.01'
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.20'
$ws.Range("E28").Value = '  +5.35%  '
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  -1.66%  '
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.87'
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("E34").Value = '  -0.65%  '
$ws.Range("D35").Value = '0.0₃0840'
$ws.Range("E35").Value = '  -5.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.29'
$ws.Range("E37").Value = '  -6.71%  '
$ws.Range("E38").Value = '  -3.08%  '
$ws.Range("E39").Value = '  -2.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '51.02'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '438.91'
$ws.Range("E41").Value = '  -1.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.79'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").Value = '2.903.08'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("E46").Value = '  -3.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.92'
$ws.Range("E47").Value = '  +3.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.99'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.20'
$ws.Range("E51").Value = '  -2.85%  '
